$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set K3 to "Done" (new shared string added for avg weekly sales column / planning fix)
$ws.Range("K3").Value = "Done"

# Update the active cell selection from E2 to I2
$ws.Range("I2").Select()
